$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

try {
  $excel.ActiveWindow.ScrollRow = 121
  $r = $excel.ActiveWindow.ScrollRow
  Write-Host "ScrollRow after set:"
  Write-Host $r
} catch {
  Write-Host "ScrollRow set failed: $_"
}

try {
  $excel.ActiveWindow.TopLeftCell = "A121"
  Write-Host "set TopLeftCell as string ok"
} catch {
  Write-Host "TopLeftCell string set failed: $_"
}
